$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.008.15'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '3.483.89'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''413.73'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").Value = '''130.15'
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = '''0.626'
$ws.Range("E7").Value = '  -2.78%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").Value = '''0.154'
$ws.Range("E10").Value = '  +7.68%  '
$ws.Range("D11").Value = '''42.49'
$ws.Range("E11").Value = '  -4.36%  '
$ws.Range("D12").Value = '''9.78'
$ws.Range("E12").Value = '  +3.79%  '
$ws.Range("D13").Value = '''0.0000225'
$ws.Range("E13").Value = '  -0.88%  '
$ws.Range("D14").Value = '4.036.17'
$ws.Range("E14").Value = '  +1.83%  '
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '''20.47'
$ws.Range("E16").Value = '  -4.39%  '
$ws.Range("D17").Value = '3.482.30'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '''12.61'
$ws.Range("E18").Value = '  +0.20%  '
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '62.927.68'
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").Value = '''468.64'
$ws.Range("E21").Value = '  -6.52%  '
$ws.Range("D22").Value = '''90.55'
$ws.Range("E22").Value = '  -3.98%  '
$ws.Range("E23").Value = '  +2.59%  '
$ws.Range("D24").Value = '''13.14'
$ws.Range("E24").Value = '  -0.97%  '
$ws.Range("E25").Value = '  +12.53%  '
$ws.Range("D26").Value = '''3.31'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("E28").Value = '  +0.35%  '
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("E30").Value = '  -0.34%  '
$ws.Range("D31").Value = '''2.68'
$ws.Range("E31").Value = '  -1.67%  '
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("E33").Value = '  -2.29%  '
$ws.Range("D34").Value = '''40.69'
$ws.Range("E34").Value = '  -4.98%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").Value = '''58.06'
$ws.Range("E37").Value = '  -4.76%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '''2.81'
$ws.Range("E39").Value = '  +7.84%  '
$ws.Range("E40").Value = '  +2.38%  '
$ws.Range("D41").Value = '''150.74'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("D43").Value = '''0.319'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("E44").Value = '  -3.41%  '
$ws.Range("D45").Value = '''4.41'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").Value = '0.0₃0588'
$ws.Range("E47").Value = '  +34.09%  '
$ws.Range("E48").Value = '  +11.44%  '
$ws.Range("D49").Value = '''16.38'
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("D50").Value = '''22.25'
$ws.Range("E50").Value = '  -1.77%  '
$ws.Range("E51").Value = '  -3.14%  '
